$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New section: uu/nhuoc diem of VBScript + Windows Script automation on IE
$ws.Range("D28").Value = 'Dễ học và sử dụng: VBScript là ngôn ngữ đơn giản và dễ học, đặc biệt là đối với những người làm việc trong môi trường Windows.'
$ws.Range("D29").Value = 'Sẵn có trong Windows: VBScript và Windows Script Host (WSH) thường đã được cài đặt sẵn trên hầu hết các máy tính Windows, do đó không cần phải cài đặt thêm phần mềm.'
$ws.Range("D30").Value = 'Tích hợp hệ thống: VBScript có thể dễ dàng tích hợp vào các tác vụ hệ thống Windows khác nhau như quản lý tệp, đối tượng Active Directory và các dịch vụ hệ thống khác.'
$ws.Range("D31").Value = 'Thao tác với IE COM object: Windows Script Host cho phép bạn tương tác trực tiếp với trình duyệt IE thông qua COM object, giúp bạn thực hiện các tác vụ trên trình duyệt.'
$ws.Range("A27").Value = 1
$ws.Range("B27").Value = 'Ưu điểm và nhược điểm khi tự động hóa trên trình duyệt IE với Vbscript và Windows Script'
$ws.Range("D34").Value = 'Khả năng tương thích hạn chế: Trình duyệt IE không còn được phát triển và không hỗ trợ nhiều tiêu chuẩn web mới như các trình duyệt hiện đại khác. Điều này có thể làm cho mã VBScript và Windows Script không tương thích với các trang web mới.'
$ws.Range("D35").Value = 'Khả năng tương thích với trình duyệt khác: VBScript và Windows Script chủ yếu được thiết kế cho IE, không phải là một giải pháp chung cho việc tự động hóa trình duyệt, đặc biệt khi bạn cần tự động hóa trình duyệt khác như Chrome hoặc Firefox.'
$ws.Range("D36").Value = 'Khả năng mở rộng hạn chế: VBScript là một ngôn ngữ đơn giản và không có các tính năng phức tạp như các ngôn ngữ lập trình hiện đại khác, giới hạn trong việc xây dựng ứng dụng phức tạp.'
$ws.Range("D37").Value = 'Khả năng hiển thị và gỡ lỗi hạn chế: Windows Script Host không cung cấp môi trường phát triển tương tự như các IDE hiện đại, dẫn đến khó khăn trong việc hiển thị mã và gỡ lỗi.'
$ws.Range("D39").Value = 'tm lại, việc sử dụng VBScript và Windows Script để tự động hóa trình duyệt IE có thể phù hợp cho các tác vụ đơn giản trong môi trường Windows truyền thống. Tuy nhiên, để thực hiện tự động hóa trên các trình duyệt khác và có tính mở rộng cao hơn, bạn nên xem xét sử dụng các ngôn ngữ lập trình hiện đại hơn và các công cụ như Selenium WebDriver.'

$ws.Range("D40").Select()

